$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "595.53"). Excel would
# normally auto-convert such text typed into a General-formatted cell into a
# real number, which would change the cell type away from text and could also
# lose exact formatting (trailing zeros, leading "0.", etc). Mark these cells as
# Text first so the assigned value is kept verbatim as a string, matching the
# original inline-string cell type, then restore the default "Normal" style so
# no visible/semantic formatting changes are introduced.
$textCells = @(
    @(5, 4),
    @(6, 4),
    @(7, 4),
    @(8, 4),
    @(9, 4),
    @(10, 4),
    @(11, 4),
    @(13, 4),
    @(15, 4),
    @(18, 4),
    @(19, 4),
    @(20, 4),
    @(21, 4),
    @(23, 4),
    @(24, 4),
    @(25, 4),
    @(26, 4),
    @(27, 4),
    @(30, 4),
    @(33, 4),
    @(34, 4),
    @(35, 4),
    @(36, 4),
    @(37, 4),
    @(38, 4),
    @(39, 4),
    @(40, 4),
    @(42, 4),
    @(43, 4),
    @(44, 4),
    @(45, 4),
    @(46, 4),
    @(47, 4),
    @(48, 4),
    @(49, 4),
    @(50, 4),
    @(51, 4)
)
foreach ($coord in $textCells) {
    $ws.Cells.Item($coord[0], $coord[1]).NumberFormat = "@"
}

# Row 2
$ws.Cells.Item(2, 4).Value = "64.622.41"
$ws.Cells.Item(2, 5).Value = "  +1.07%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.622.06"
$ws.Cells.Item(3, 5).Value = "  -0.13%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "595.53"
$ws.Cells.Item(5, 5).Value = "  -0.34%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "154.59"
$ws.Cells.Item(6, 5).Value = "  +1.62%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.999"
$ws.Cells.Item(7, 5).Value = "  -0.07%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.589"
$ws.Cells.Item(8, 5).Value = "  -0.19%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.116"
$ws.Cells.Item(9, 5).Value = "  +5.79%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "0.397"
$ws.Cells.Item(10, 5).Value = "  +2.80%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "5.78"
$ws.Cells.Item(11, 5).Value = "  +1.28%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +1.18%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "28.73"
$ws.Cells.Item(13, 5).Value = "  +2.63%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "3.090.24"
$ws.Cells.Item(14, 5).Value = "  -0.25%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "ShibaInu"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(15, 4).Value = "0.0000172"
$ws.Cells.Item(15, 5).Value = "  +11.34%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "WrappedBTC"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(16, 4).Value = "64.490.13"
$ws.Cells.Item(16, 5).Value = "  +1.13%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.590.30"
$ws.Cells.Item(17, 5).Value = "  -0.89%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "12.44"
$ws.Cells.Item(18, 5).Value = "  -0.06%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "4.80"
$ws.Cells.Item(19, 5).Value = "  +1.31%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "351.69"
$ws.Cells.Item(20, 5).Value = "  +0.83%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "7.16"
$ws.Cells.Item(21, 5).Value = "  +3.48%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +0.16%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "67.58"
$ws.Cells.Item(23, 5).Value = "  +0.35%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "1.72"
$ws.Cells.Item(24, 5).Value = "  -0.90%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "9.42"
$ws.Cells.Item(25, 5).Value = "  -0.20%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "1.65"
$ws.Cells.Item(26, 5).Value = "  -3.44%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "8.29"
$ws.Cells.Item(27, 5).Value = "  +2.17%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +1.33%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "528.86"
$ws.Cells.Item(30, 5).Value = "  -5.46%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "0.0₃0913"
$ws.Cells.Item(31, 5).Value = "  +6.40%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -0.13%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "1.81"
$ws.Cells.Item(33, 5).Value = "  +2.85%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "5.77"
$ws.Cells.Item(34, 5).Value = "  +8.32%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "6.31"
$ws.Cells.Item(35, 5).Value = "  +1.24%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "0.425"
$ws.Cells.Item(36, 5).Value = "  +1.90%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "164.77"
$ws.Cells.Item(37, 5).Value = "  -1.39%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "Stacks"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(38, 4).Value = "2.02"
$ws.Cells.Item(38, 5).Value = "  +3.35%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "EthereumClassic"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(39, 4).Value = "20.14"
$ws.Cells.Item(39, 5).Value = "  +2.52%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "0.999"
$ws.Cells.Item(40, 5).Value = "  -0.12%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -0.03%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "165.96"
$ws.Cells.Item(42, 5).Value = "  -1.17%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "41.49"
$ws.Cells.Item(43, 5).Value = "  +4.11%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "4.07"
$ws.Cells.Item(44, 5).Value = "  +2.00%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Hedera"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(45, 4).Value = "0.0603"
$ws.Cells.Item(45, 5).Value = "  +1.84%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(46, 4).Value = "23.17"
$ws.Cells.Item(46, 5).Value = "  +5.05%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "2.23"
$ws.Cells.Item(47, 5).Value = "  +7.17%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "0.640"
$ws.Cells.Item(48, 5).Value = "  +0.85%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "0.0250"
$ws.Cells.Item(49, 5).Value = "  -0.97%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "0.0981"
$ws.Cells.Item(50, 5).Value = "  +1.02%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "19.30"
$ws.Cells.Item(51, 5).Value = "  -0.91%  "

# Restore default styling on the cells we temporarily marked as Text
foreach ($coord in $textCells) {
    $ws.Cells.Item($coord[0], $coord[1]).Style = "Normal"
}
